$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Count number of occurrences (or frequency) in a sorted array"
$ws.Range("B8").Value = "CountNumberOfOccurrences"

$ws.Range("A5").Select()
